# Add Sheet3 (after Sheet2): a new table of ReactOS module paths together
# with the specific short-name resource/translation files that now get
# run through `lang_add`/`translate2`, so every language file translation
# is driven from inside the spreadsheet tool.
#
# Cell values are written column-by-column (A, then B+C per row, then D+E
# per row) to match the order the rows were originally authored in, which
# also determines the order new entries land in the shared-string table.

$wb = $excel.ActiveWorkbook

$data = @(
    @('C:\sources\reactos\base\applications\winhlp32','rsrc.rc','En.rc','cd C:\sources\reactos\base\applications\winhlp32','call :lang_add rsrc.rc En.rc'),
    @('C:\sources\reactos\base\applications\write\','rsrc.rc','En.rc','cd C:\sources\reactos\base\applications\write\','call :lang_add rsrc.rc En.rc'),
    @('C:\sources\reactos\base\setup\usetup\','usetup.rc','en-US.h','cd C:\sources\reactos\base\setup\usetup\','call :lang_add usetup.rc en-US.h'),
    @('C:\sources\reactos\dll\win32\avifil32\','rsrc.rc','avifile_En.rc','cd C:\sources\reactos\dll\win32\avifil32\','call :lang_add rsrc.rc avifile_En.rc'),
    @('C:\sources\reactos\dll\win32\comctl32\','rsrc.rc','comctl_En.rc','cd C:\sources\reactos\dll\win32\comctl32\','call :lang_add rsrc.rc comctl_En.rc'),
    @('C:\sources\reactos\dll\win32\comdlg32\','rsrc.rc','cdlg_En.rc','cd C:\sources\reactos\dll\win32\comdlg32\','call :lang_add rsrc.rc cdlg_En.rc'),
    @('C:\sources\reactos\dll\win32\credui\','credui.rc','credui_En.rc','cd C:\sources\reactos\dll\win32\credui\','call :lang_add credui.rc credui_En.rc'),
    @('C:\sources\reactos\dll\win32\crypt32\','crypt32.rc','crypt32_En.rc','cd C:\sources\reactos\dll\win32\crypt32\','call :lang_add crypt32.rc crypt32_En.rc'),
    @('C:\sources\reactos\dll\win32\cryptdlg\','cryptdlg.rc','cryptdlg_En.rc','cd C:\sources\reactos\dll\win32\cryptdlg\','call :lang_add cryptdlg.rc cryptdlg_En.rc'),
    @('C:\sources\reactos\dll\win32\cryptui\','cryptui.rc','cryptui_En.rc','cd C:\sources\reactos\dll\win32\cryptui\','call :lang_add cryptui.rc cryptui_En.rc'),
    @('C:\sources\reactos\dll\win32\hhctrl.ocx\','hhctrl.rc','En.rc','cd C:\sources\reactos\dll\win32\hhctrl.ocx\','call :lang_add hhctrl.rc En.rc'),
    @('C:\sources\reactos\dll\win32\iccvid\','rsrc.rc','iccvid_En.rc','cd C:\sources\reactos\dll\win32\iccvid\','call :lang_add rsrc.rc iccvid_En.rc'),
    @('C:\sources\reactos\dll\win32\jscript\','rsrc.rc','jscript_En.rc','cd C:\sources\reactos\dll\win32\jscript\','call :lang_add rsrc.rc jscript_En.rc'),
    @('C:\sources\reactos\dll\win32\mapi32\','version.rc','En.rc','cd C:\sources\reactos\dll\win32\mapi32\','call :lang_add version.rc En.rc'),
    @('C:\sources\reactos\dll\win32\mpr\','version.rc','mpr_En.rc','cd C:\sources\reactos\dll\win32\mpr\','call :lang_add version.rc mpr_En.rc'),
    @('C:\sources\reactos\dll\win32\msacm32\','msacm.rc','msacm_En.rc','cd C:\sources\reactos\dll\win32\msacm32\','call :lang_add msacm.rc msacm_En.rc'),
    @('C:\sources\reactos\dll\win32\mshtml\','rsrc.rc','En.rc','cd C:\sources\reactos\dll\win32\mshtml\','call :lang_add rsrc.rc En.rc'),
    @('C:\sources\reactos\dll\win32\msi\','msi.rc','msi_En.rc','cd C:\sources\reactos\dll\win32\msi\','call :lang_add msi.rc msi_En.rc'),
    @('C:\sources\reactos\dll\win32\msrle32\','rsrc.rc','msrle_En.rc','cd C:\sources\reactos\dll\win32\msrle32\','call :lang_add rsrc.rc msrle_En.rc'),
    @('C:\sources\reactos\dll\win32\msvfw32\','rsrc.rc','msvfw32_En.rc','cd C:\sources\reactos\dll\win32\msvfw32\','call :lang_add rsrc.rc msvfw32_En.rc'),
    @('C:\sources\reactos\dll\win32\msvidc32\','rsrc.rc','msvidc32_En.rc','cd C:\sources\reactos\dll\win32\msvidc32\','call :lang_add rsrc.rc msvidc32_En.rc'),
    @('C:\sources\reactos\dll\win32\oleacc\','oleacc.rc','oleacc_En.rc','cd C:\sources\reactos\dll\win32\oleacc\','call :lang_add oleacc.rc oleacc_En.rc'),
    @('C:\sources\reactos\dll\win32\oleaut32\','oleaut32.rc','oleaut32_En.rc','cd C:\sources\reactos\dll\win32\oleaut32\','call :lang_add oleaut32.rc oleaut32_En.rc'),
    @('C:\sources\reactos\dll\win32\oledlg\','rsrc.rc','oledlg_En.rc','cd C:\sources\reactos\dll\win32\oledlg\','call :lang_add rsrc.rc oledlg_En.rc'),
    @('C:\sources\reactos\dll\win32\shdoclc\','rsrc.rc','En.rc','cd C:\sources\reactos\dll\win32\shdoclc\','call :lang_add rsrc.rc En.rc'),
    @('C:\sources\reactos\dll\win32\shlwapi\','version.rc','shlwapi_En.rc','cd C:\sources\reactos\dll\win32\shlwapi\','call :lang_add version.rc shlwapi_En.rc'),
    @('C:\sources\reactos\dll\win32\wininet\','rsrc.rc','wininet_En.rc','cd C:\sources\reactos\dll\win32\wininet\','call :lang_add rsrc.rc wininet_En.rc'),
    @('C:\sources\reactos\dll\win32\winmm\','winmm_res.rc','winmm_En.rc','cd C:\sources\reactos\dll\win32\winmm\','call :lang_add winmm_res.rc winmm_En.rc'),
    @('C:\sources\reactos\dll\win32\wldap32\','wldap32.rc','wldap32_En.rc','cd C:\sources\reactos\dll\win32\wldap32\','call :lang_add wldap32.rc wldap32_En.rc'),
    @('C:\sources\reactos\win32ss\printing\monitors\localmon\ui\','localui.rc','ui_En.rc','cd C:\sources\reactos\win32ss\printing\monitors\localmon\ui\','call :lang_add localui.rc ui_En.rc')
)


$afterSheet = $wb.Worksheets.Item(2)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $afterSheet)
$ws.Name = "Sheet3"

for ($i = 0; $i -lt $data.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $data[$i][0]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $data[$i][1]
    $ws.Cells.Item($i + 1, 3).Value = $data[$i][2]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $ws.Cells.Item($i + 1, 4).Value = $data[$i][3]
    $ws.Cells.Item($i + 1, 5).Value = $data[$i][4]
}

$ws.Columns("A:C").AutoFit()
$ws.Range("C15").Select()
$ws.Activate()
